# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with refreshed values. Price cells are forced to Text format
# before assignment (and reset to the Normal style afterwards) so that
# values such as "242.50" or "0.9995" keep their exact textual
# representation instead of being coerced into numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '29.920.58'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.07%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.875.92'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("E4").Value = '  +0.04%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.7435'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -3.86%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '242.50'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3150'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +1.10%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.07238'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("E10").Value = '  -3.73%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.08386'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -2.69%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.7524'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -1.56%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '5.426'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +0.93%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '1.872.91'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -3.35%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '92.51'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -1.42%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '29.914.42'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -0.22%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '6.081'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -1.68%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '249.33'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +1.90%  '
$ws.Range("E19").Value = '  -1.41%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.000007857'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +0.40%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '0.9995'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +0.11%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '2.127.59'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -3.98%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '8.037'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +0.22%  '
$ws.Range("E24").Value = '  -0.03%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.1561'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -5.34%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '9.264'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -1.24%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '165.09'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +1.90%  '
$ws.Range("E28").Value = '  -0.67%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '2.035'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +0.04%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '1.517'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +5.10%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '4.598'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +1.51%  '
$ws.Range("E32").Value = '  +0.20%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '4.283'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +4.38%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.05332'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -1.75%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.237'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -0.30%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.7506'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +0.66%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.9992'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("E38").Value = '  -0.24%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.01968'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -0.01%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '2.759'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.81%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.4538'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +1.59%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '1.114.07'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +0.29%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '6.055'
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '72.45'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -1.22%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.8552'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +0.44%  '
$ws.Range("E46").Value = '  +0.10%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '103.33'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.46%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.857'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -0.79%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '7.620'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -0.13%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '2.024.80'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -3.44%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '2.902'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -2.75%  '
